# Apply "Updates with wiring updates" commit to the "pins" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pins")
$ws.Activate()

# Row 4 (pin 2 / PTD0): clear the old FLTSD1 / FLT-SD/ISO1 wiring info
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Row 5 (pin 3 / PTA12): HIN1 -> 2INA, Hin1_ISO -> toms driver
$ws.Range("D5").Value = "2INA"
$ws.Range("E5").Value = "toms driver  "

# Row 6 (pin 4 / PTA13): LIN1 -> 2INB, Lin_ISO -> toms driver
$ws.Range("D6").Value = "2INB"
$ws.Range("E6").Value = "toms driver  "

# Row 7 (pin 5 / PTD7): FLTCLR1 -> 2PWM, FLT_CLR1 -> toms driver
$ws.Range("D7").Value = "2PWM"
$ws.Range("E7").Value = "toms driver  "

# Row 9 (pin 7 / PTD2): clear the old FLTSD2 / FLT-SD/ISO2 wiring info
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()

# Row 10 (pin 8 / PTD3): HIN2 -> 1INA, Hin2_ISO -> toms driver
$ws.Range("D10").Value = "1INA"
$ws.Range("E10").Value = "toms driver  "

# Row 11 (pin 9 / PTC3): LIN2 -> 1PWM, Lin2_ISO -> toms driver
$ws.Range("D11").Value = "1PWM"
$ws.Range("E11").Value = "toms driver  "

# Row 12 (pin 10 / PTC4): FLTCLR2 -> 1INB, FLT_CLR2 -> toms driver
$ws.Range("D12").Value = "1INB"
$ws.Range("E12").Value = "toms driver  "

# Row 18 (pin 16 / PTB0): TEMP1 -> TEMP2
$ws.Range("D18").Value = "TEMP2"

# Row 21 (pin 19 / PTB2): TEMP2 -> TEMP1
$ws.Range("D21").Value = "TEMP1"

# Update the selected cell on the pins sheet
$ws.Range("B18").Select()
